$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns stay text (avoid Excel auto-numeric conversion)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "63.199.50"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").Value = "3.161.36"
$ws.Range("E3").Value = "  -1.22%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "597.31"
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("D6").Value = "134.93"
$ws.Range("E6").Value = "  -1.12%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "3.158.78"
$ws.Range("E8").Value = "  -1.19%  "
$ws.Range("E9").Value = "  +1.42%  "
$ws.Range("E10").Value = "  -1.72%  "
$ws.Range("E11").Value = "  -0.63%  "
$ws.Range("D12").Value = "0.453"
$ws.Range("E12").Value = "  -0.45%  "
$ws.Range("D13").Value = "0.0000238"
$ws.Range("E13").Value = "  -0.01%  "
$ws.Range("D14").Value = "34.50"
$ws.Range("E14").Value = "  +2.61%  "
$ws.Range("D15").Value = "3.679.78"
$ws.Range("E15").Value = "  -1.12%  "
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("D17").Value = "3.164.22"
$ws.Range("E17").Value = "  -0.69%  "
$ws.Range("D18").Value = "63.177.89"
$ws.Range("E18").Value = "  +1.03%  "
$ws.Range("E19").Value = "  -2.13%  "
$ws.Range("D20").Value = "459.51"
$ws.Range("E20").Value = "  -0.72%  "
$ws.Range("D21").Value = "13.91"
$ws.Range("E21").Value = "  -0.84%  "
$ws.Range("E22").Value = "  -2.78%  "
$ws.Range("D23").Value = "7.62"
$ws.Range("E23").Value = "  -1.06%  "
$ws.Range("D24").Value = "82.95"
$ws.Range("E24").Value = "  -1.04%  "
$ws.Range("D25").Value = "13.19"
$ws.Range("E25").Value = "  -2.58%  "
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("E27").Value = "  -1.07%  "
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("E29").Value = "  -2.71%  "
$ws.Range("E30").Value = "  -0.87%  "
$ws.Range("D31").Value = "7.69"
$ws.Range("E31").Value = "  -3.56%  "
$ws.Range("D32").Value = "27.10"
$ws.Range("E32").Value = "  -1.17%  "
$ws.Range("E33").Value = "  -2.23%  "
$ws.Range("D34").Value = "2.40"
$ws.Range("E34").Value = "  -2.15%  "
$ws.Range("D36").Value = "5.86"
$ws.Range("E36").Value = "  -0.33%  "
$ws.Range("E37").Value = "  +5.63%  "
$ws.Range("D38").Value = "51.22"
$ws.Range("E38").Value = "  -0.85%  "
$ws.Range("D39").Value = "0.0388"
$ws.Range("E39").Value = "  -0.60%  "
$ws.Range("D40").Value = "8.12"
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("E41").Value = "  -2.66%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").Value = "391.21"
$ws.Range("E42").Value = "  -6.04%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "2.61"
$ws.Range("E43").Value = "  -1.84%  "
$ws.Range("D44").Value = "2.788.16"
$ws.Range("E44").Value = "  -7.14%  "
$ws.Range("E45").Value = "  -1.24%  "
$ws.Range("D46").Value = "127.21"
$ws.Range("E46").Value = "  +2.33%  "
$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D47").Value = "35.79"
$ws.Range("E47").Value = "  -1.34%  "
$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").Value = "0.999"
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("D49").Value = "2.09"
$ws.Range("E49").Value = "  -4.38%  "
$ws.Range("E50").Value = "  -0.67%  "
$ws.Range("E51").Value = "  -4.19%  "
